$wb = $excel.ActiveWorkbook

$oldGuid = "a4a4156c-90a4-405b-ae8d-fd882bce9998"
$newGuid = "d169085d-ee20-480a-8fed-8c9db05c8fc5"

$oldHash = "0050d138b3604575d523da5cd32743a6c1421f7c"
$newHash = "2e8e0b8f07559529eb2e026432081d639ed7e6ef"

# The external hyperlink target (commit SHA + path) itself is not touched by
# the change - only the cell text / hyperlink display text gets the new guid.
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e068c132258c8608bb8c6cc35c8981473604c946/e2e/$oldGuid.md"

$sOverview = $wb.Worksheets.Item("Overview")
$sZhCn = $wb.Worksheets.Item("zh-cn")
$sDeDe = $wb.Worksheets.Item("de-de")

# ---- Overview sheet ----
# A2 = guid.md ; B2 = e2e\guid.md (hyperlink) ; G2 = handoff datetime
$sOverview.Range("A2").Value = "$newGuid.md"
$sOverview.Range("B2").Value = "e2e\$newGuid.md"
$sOverview.Hyperlinks.Delete()
$sOverview.Hyperlinks.Add($sOverview.Range("B2"), $hyperlinkAddress, "", "", "e2e\$newGuid.md")
$sOverview.Range("G2").Value = "2016-08-27 14:56:14"

# ---- zh-cn sheet ----
# A2 = guid.md (hyperlink) ; G2 = handoff xlf file ; H2 = handoff datetime
$sZhCn.Range("A2").Value = "$newGuid.md"
$sZhCn.Hyperlinks.Delete()
$sZhCn.Hyperlinks.Add($sZhCn.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")
$sZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$sZhCn.Range("H2").Value = "2016-08-27 14:56:10"

# ---- de-de sheet ----
# A2 = guid.md (hyperlink) ; G2 = handoff xlf file ; H2 = handoff datetime
$sDeDe.Range("A2").Value = "$newGuid.md"
$sDeDe.Hyperlinks.Delete()
$sDeDe.Hyperlinks.Add($sDeDe.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")
$sDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$sDeDe.Range("H2").Value = "2016-08-27 14:56:14"
